$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.913.65"
$ws.Range("E2").Value = "  +2.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.75"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.85"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.79"
$ws.Range("E8").Value = "  +4.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("E9").Value = "  +0.38%  "

$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.903.38"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.665.75"
$ws.Range("E13").Value = "  -1.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.15"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.14"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "253.24"
$ws.Range("E17").Value = "  +7.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.866.64"
$ws.Range("E18").Value = "  +2.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0732"
$ws.Range("E19").Value = "  -1.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  -3.53%  "

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("E22").Value = "  -1.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("E24").Value = "  -1.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.90"
$ws.Range("E25").Value = "  -1.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.24"
$ws.Range("E26").Value = "  -2.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -0.66%  "

$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("E30").Value = "  +5.93%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("E32").Value = "  -0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.15"
$ws.Range("E33").Value = "  -2.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.420.80"
$ws.Range("E34").Value = "  -7.77%  "

$ws.Range("E35").Value = "  -4.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.933"
$ws.Range("E37").Value = "  -0.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.581"
$ws.Range("E38").Value = "  -4.45%  "

$ws.Range("E39").Value = "  -1.22%  "

$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "69.67"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("E41").Value = "  -2.97%  "

$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("E43").Value = "  -0.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.811.61"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.38"
$ws.Range("E45").Value = "  -6.79%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.791"
$ws.Range("E46").Value = "  +1.52%  "

$ws.Range("E47").Value = "  +4.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.96"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("E49").Value = "  -0.51%  "

$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0511"
$ws.Range("E51").Value = "  +0.10%  "
